$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 2.98
$ws.Range("G3").Value = 4
$ws.Range("I3").Value = 2.68
$ws.Range("P3").Value = 1.75
$ws.Range("Q3").Value = 1.95
$ws.Range("V3").Value = 1.59
$ws.Range("W3").Value = 1.33
$ws.Range("G5").Value = 2.7
$ws.Range("H5").Value = 2.56
$ws.Range("I5").Value = 3.1
$ws.Range("J5").Value = 3.95
$ws.Range("V5").Value = 1.47
$ws.Range("W5").Value = 1.58
$ws.Range("F7").Value = 5.7
$ws.Range("G7").Value = 38
$ws.Range("H7").Value = 1.34
$ws.Range("I7").Value = 1.45
$ws.Range("P7").Value = 1.89
$ws.Range("Q7").Value = 1.89
$ws.Range("S7").Value = 1.9
$ws.Range("P9").Value = 2.2
$ws.Range("Q9").Value = 1.7
$ws.Range("AA9").Value = 60
$ws.Range("AO9").Value = 29
$ws.Range("F10").Value = 1.84
$ws.Range("G10").Value = 1.86
$ws.Range("J10").Value = 3.6
$ws.Range("K10").Value = 3.85
$ws.Range("W10").Value = 2.16
$ws.Range("G11").Value = 1.93
$ws.Range("L11").Value = 1.41
$ws.Range("Q11").Value = 2.02
$ws.Range("R11").Value = 1.3
$ws.Range("V11").Value = 1.18
$ws.Range("R14").Value = 1.53
$ws.Range("T14").Value = 1.43
$ws.Range("F15").Value = 3.3
$ws.Range("G15").Value = 3.95
$ws.Range("H15").Value = 1.93
$ws.Range("I15").Value = 2.14
$ws.Range("K15").Value = 5.6
$ws.Range("N15").Value = 3.6
$ws.Range("O15").Value = 1.08
$ws.Range("P15").Value = 3.6
$ws.Range("R15").Value = 2.02
$ws.Range("V15").Value = 1.88
$ws.Range("W15").Value = 1.37
$ws.Range("X15").Value = 70
$ws.Range("Y15").Value = 980
$ws.Range("Z15").Value = 980
$ws.Range("AA15").Value = 980
$ws.Range("AB15").Value = 980
$ws.Range("AC15").Value = 980
$ws.Range("AD15").Value = 980
$ws.Range("AE15").Value = 980
$ws.Range("AF15").Value = 55
$ws.Range("AG15").Value = 980
$ws.Range("AH15").Value = 980
$ws.Range("AI15").Value = 980
$ws.Range("AJ15").Value = 90
$ws.Range("AK15").Value = 980
$ws.Range("AL15").Value = 980
$ws.Range("AM15").Value = 55
$ws.Range("G16").Value = 4.9
$ws.Range("I16").Value = 2.86
$ws.Range("J16").Value = 2.72
$ws.Range("K16").Value = 980
$ws.Range("N16").Value = 1.35
$ws.Range("O16").Value = 1.01
$ws.Range("Q16").Value = 2.42
$ws.Range("S16").Value = 2.42
$ws.Range("V16").Value = 1.54
$ws.Range("W16").Value = 1.25
$ws.Range("F17").Value = 1.92
$ws.Range("K17").Value = 980
$ws.Range("N17").Value = 1.58
$ws.Range("P17").Value = 1.58
$ws.Range("Q17").Value = 1.99
$ws.Range("S17").Value = 2
$ws.Range("F18").Value = 1.5
$ws.Range("I18").Value = 9
$ws.Range("L18").Value = 1.43
$ws.Range("T18").Value = 2.34
$ws.Range("U18").Value = 1.69
$ws.Range("V18").Value = 1.12
$ws.Range("W18").Value = 2.96
$ws.Range("X18").Value = 13
$ws.Range("Z18").Value = 75
$ws.Range("AA18").Value = 400
$ws.Range("AB18").Value = 6.8
$ws.Range("AC18").Value = 9.800000000000001
$ws.Range("AD18").Value = 34
$ws.Range("AE18").Value = 190
$ws.Range("AH18").Value = 32
$ws.Range("AI18").Value = 170
$ws.Range("AJ18").Value = 12
$ws.Range("AK18").Value = 18.5
$ws.Range("AL18").Value = 48
$ws.Range("AM18").Value = 230
$ws.Range("AO18").Value = 300
$ws.Range("F20").Value = 2.06
$ws.Range("G20").Value = 2.08
$ws.Range("I20").Value = 4.5
$ws.Range("K20").Value = 3.5
$ws.Range("L20").Value = 1.01
$ws.Range("N20").Value = 3.3
$ws.Range("R20").Value = 1.29
$ws.Range("V20").Value = 1.28
$ws.Range("W20").Value = 1.93
$ws.Range("H21").Value = 1.52
$ws.Range("I21").Value = 1.61
$ws.Range("L21").Value = 1.01
$ws.Range("M21").Value = 1.06
$ws.Range("N21").Value = 3.85
$ws.Range("O21").Value = 1.28
$ws.Range("Q21").Value = 1.83
$ws.Range("R21").Value = 1.39
$ws.Range("S21").Value = 3.1
$ws.Range("T21").Value = 1.9
$ws.Range("U21").Value = 1.89
$ws.Range("V21").Value = 2.62
$ws.Range("W21").Value = 1.14
$ws.Range("X21").Value = 20
$ws.Range("Y21").Value = 8.6
$ws.Range("Z21").Value = 9.6
$ws.Range("AA21").Value = 15.5
$ws.Range("AB21").Value = 28
$ws.Range("AC21").Value = 10.5
$ws.Range("AD21").Value = 10.5
$ws.Range("AE21").Value = 980
$ws.Range("AF21").Value = 75
$ws.Range("AG21").Value = 980
$ws.Range("AH21").Value = 980
$ws.Range("AI21").Value = 980
$ws.Range("AJ21").Value = 270
$ws.Range("AK21").Value = 140
$ws.Range("AL21").Value = 120
$ws.Range("AM21").Value = 160
$ws.Range("AN21").Value = 170
$ws.Range("AO21").Value = 8.800000000000001
$ws.Range("L22").Value = 1.39
$ws.Range("W22").Value = 1.88
$ws.Range("AD22").Value = 20
$ws.Range("AN22").Value = 17
